$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Student ID values (column A, rows 2-7) ---
# Force text number-format first so the numeric-looking IDs are stored
# as text (matching the original inline-string typing) instead of
# being auto-converted to numbers.
$ws.Range("A2:A7").NumberFormat = "@"

$ws.Range("A2").Value = "200869"
$ws.Range("A3").Value = "200852"
$ws.Range("A4").Value = "200897"
$ws.Range("A5").Value = "200850"
$ws.Range("A6").Value = "211137"
$ws.Range("A7").Value = "201838"

# --- Append 3 new rows (8, 9, 10) with the same layout/formatting ---
# Row 8 and 10 mirror the styling of the existing even rows (e.g. row 6),
# row 9 mirrors the styling of the existing odd rows (e.g. row 7).
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A8:F8").PasteSpecial(-4122) | Out-Null

$ws.Range("A7:F7").Copy() | Out-Null
$ws.Range("A9:F9").PasteSpecial(-4122) | Out-Null

$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A10:F10").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A10").NumberFormat = "@"

# Row 8
$ws.Range("A8").Value = "200785"
$ws.Range("B8").Value = "general surgery"
$ws.Range("C8").Value = "13/10/2025"
$ws.Range("D8").Value = "10:30:00"
$ws.Range("E8").Value = "Excuse"
$ws.Range("F8").Value = "System"

# Row 9
$ws.Range("A9").Value = "201574"
$ws.Range("B9").Value = "general surgery"
$ws.Range("C9").Value = "13/10/2025"
$ws.Range("D9").Value = "10:30:00"
$ws.Range("E9").Value = "Excuse"
$ws.Range("F9").Value = "System"

# Row 10
$ws.Range("A10").Value = "201252"
$ws.Range("B10").Value = "general surgery"
$ws.Range("C10").Value = "13/10/2025"
$ws.Range("D10").Value = "10:30:00"
$ws.Range("E10").Value = "Excuse"
$ws.Range("F10").Value = "System"
